$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 360
$ws1.Range("F4").Value = 456
$ws1.Range("F5").Value = 1789
$ws1.Range("F7").Value = 2252
$ws1.Range("F11").Value = 5089
$ws1.Range("F12").Value = 375
$ws1.Range("F17").Value = 207
$ws1.Range("F21").Value = 4125
$ws1.Range("F22").Value = 738
$ws1.Range("F23").Value = 741
$ws1.Range("F24").Value = 36
$ws1.Range("F26").Value = 117
$ws1.Range("F27").Value = 135
$ws1.Range("F30").Value = 102
$ws1.Range("F33").Value = 29
$ws1.Range("F34").Value = 1049
$ws1.Range("F36").Value = 2664
$ws1.Range("F38").Value = 49

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 360
$ws4.Range("F4").Value = 456
$ws4.Range("F5").Value = 1789
$ws4.Range("F7").Value = 2252
$ws4.Range("F11").Value = 5089
$ws4.Range("F12").Value = 375
$ws4.Range("F17").Value = 207
$ws4.Range("F21").Value = 4126
$ws4.Range("F22").Value = 738
$ws4.Range("F23").Value = 741
$ws4.Range("F24").Value = 36
$ws4.Range("F26").Value = 117
$ws4.Range("F27").Value = 135
$ws4.Range("F30").Value = 102
$ws4.Range("F34").Value = 29
$ws4.Range("F35").Value = 1049
$ws4.Range("F37").Value = 2664
$ws4.Range("F39").Value = 49
